# Swap the deck's applied (slide-facing) theme colour scheme from the
# "Integral" / Red Violet palette over to the "Office Theme" palette that
# was previously only used by the Notes Master (ppt/theme/theme2.xml).
#
# PowerPoint's VBA/COM "RGB" colour integers are packed as
#   value = R + G*256 + B*65536
# i.e. the reverse byte order of the familiar "RRGGBB" hex notation, so we
# convert each target hex colour accordingly before assigning it.

function ConvertTo-VbaColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target theme colours ("Office Theme"), in DrawingML clrScheme slot order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $tcs.Item($i).RGB = ConvertTo-VbaColor $officeThemeColors[$i - 1]
}
